# Trade #45 closed at 2026-02-17 08:38:58 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" sheets with the latest
# aggregate statistics, and appends the newly closed trade (#45) as a new
# row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.54   # Current Capital
$summary.Range("B4").Value = -0.46     # Total P&L $
$summary.Range("B5").Value = -0.2      # Total P&L %
$summary.Range("B6").Value = 45        # Total Trades
$summary.Range("B7").Value = 16        # Winning Trades
$summary.Range("B9").Value = 35.56     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.54000000000001   # Capital
$status.Range("D4").Value = 45                  # Trades
$status.Range("E4").Value = -0.46                # P&L $
$status.Range("F4").Value = -0.46                # P&L %
$status.Range("G4").Value = 35.56                # Win Rate %

# ---------------------------------------------------------------------
# 3. Append new trade row (#45) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$newTradeRow = 46

function Add-TradeRow($ws, $row) {
    $ws.Range("A$row").Value = 45
    # Use a leading apostrophe so these remain literal text instead of
    # being auto-converted to Excel date/time serial values.
    $ws.Range("B$row").Value = "'2026-02-17"
    $ws.Range("C$row").Value = "'08:38:52"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 0.82
    $ws.Range("G$row").Value = 0.84
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = 2.439
    $ws.Range("J$row").Value = 0.02
    $ws.Range("K$row").Value = 99.54000000000001
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.13

    # Restore default (Normal) style on the text cells so the quote
    # prefix used above does not leave a lingering cell style behind.
    $ws.Range("B$row`:C$row").Style = "Normal"
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades $newTradeRow

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking $newTradeRow
